$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch/helper cell used to coerce a value into literal TEXT without
# Excel's "smart" input parser turning a numeric-looking / date-looking
# string into a real number or date, and WITHOUT leaving behind any new
# number-format style (unlike setting NumberFormat="@" directly on the
# target cell, which bakes a permanent style into the sheet). We compute
# TEXT(value,"0") in a far-away cell, copy it, paste-special VALUES ONLY
# into the destination (which carries no style along), then wipe the
# helper cell so it leaves no trace in the sheet's used range.
$helper = $ws.Range("ZZ1")

function Set-TextValue($addr, $val) {
    $dest = $ws.Range($addr)
    $escaped = $val -replace '"', '""'
    $helper.Formula = '=TEXT("' + $escaped + '","0")'
    $helper.Copy()
    $dest.PasteSpecial(-4163)
    $helper.ClearContents()
    $dest.ClearFormats()
}

# Helper: write a genuine numeric value.
function Set-NumberValue($addr, $val) {
    $ws.Range($addr).Value = $val
}

# --- Row 6: existing row, Pajak_Terhutang / Tanggal_Jatuh_Tempo / Pajak
#     switch from numbers to their literal text representations ---
Set-TextValue "E6" "90000"
Set-TextValue "F6" "2026-07-30 00:00:00"
Set-TextValue "G6" "90000"

# --- Row 7: new record (Sinta Maharani) -- all text ---
Set-TextValue "A7" "1245367800112234"
Set-TextValue "B7" "BG7783GI"
Set-TextValue "C7" "Sinta Maharani"
Set-TextValue "D7" "Palembang"
Set-TextValue "E7" "40000"
Set-TextValue "F7" "2026-07-31 00:00:00"
Set-TextValue "G7" "40000"

# --- Row 8: new record (Dinda) -- identity columns text, amounts numeric ---
Set-TextValue "A8" "9801234567819235"
Set-TextValue "B8" "BG8799BI"
Set-TextValue "C8" "Dinda"
Set-TextValue "D8" "Palembang"
Set-NumberValue "E8" 80000
$ws.Range("F8").NumberFormat = "YYYY-MM-DD"
Set-NumberValue "F8" 46234
Set-NumberValue "G8" 80000

$ws.Range("ZZ1").Clear()

Write-Host "Applied data_kendaraan update: rows 7-8 added, row 6 amounts/date converted to text"
